$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A265 / A266: inline-string "49" -> real number 49 ---
$ws.Range("A265").Value2 = 49
$ws.Range("A266").Value2 = 49

# --- Row 267 (new) ---
$ws.Range("A267").Value2 = '''51'
$ws.Range("B267").Value2 = 'Cafe / Food Video Editing For TikTok Videos / IG Reels  - Upwork'
$ws.Range("C267").Value2 = 'https://www.upwork.com/jobs/Cafe-Food-Video-Editing-For-TikTok-Videos-Reels_%7E01af5cfc4a0962de74?source=rss'
$v1 = @'
Hi, I am looking for a video editor for Instagram  and TikTok Food Videos. 
Video content revolves around cafes and restaurants. About 20-30 sec
The only thing is turn-around time has to be fast, and person has to be responsive.
This can be a multiple-video project, if videos can be done professionally and up to speed. 
To be paid by per video, experience with editing food video will be helpful. 
Interested applicants can refer to https://www.tiktok.com/@danielfooddiary for reference on the usual style of editing and what to expect. Thanks.
Budget
: $20
Posted On
: June 15, 2024 23:06 UTC
Category
: Video Editing
Skills
:Video Editing    
Skills
:        Video Editing            
Country
: Singapore
click to apply

'@
$ws.Range("D267").Value2 = $v1
$v2 = @'
Hi, I am looking for a video editor for Instagram&nbsp;&nbsp;and TikTok Food Videos. <br /><br />
Video content revolves around cafes and restaurants. About 20-30 sec<br /><br />
The only thing is turn-around time has to be fast, and person has to be responsive.<br /><br />
This can be a multiple-video project, if videos can be done professionally and up to speed. <br /><br />
To be paid by per video, experience with editing food video will be helpful. <br /><br />
Interested applicants can refer to https://www.tiktok.com/@danielfooddiary for reference on the usual style of editing and what to expect. Thanks.<br /><br /><b>Budget</b>: $20
<br /><b>Posted On</b>: June 15, 2024 23:06 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Video Editing    
<br /><b>Skills</b>:        Video Editing            <br /><b>Country</b>: Singapore
<br /><a href="https://www.upwork.com/jobs/Cafe-Food-Video-Editing-For-TikTok-Videos-Reels_%7E01af5cfc4a0962de74?source=rss">click to apply</a>

'@
$ws.Range("E267").Value2 = $v2
$ws.Range("F267").Value2 = 'Sat, 15 Jun 2024 23:06:17 +0000'
$ws.Range("G267").Value2 = 'https://www.upwork.com/jobs/Cafe-Food-Video-Editing-For-TikTok-Videos-Reels_%7E01af5cfc4a0962de74?source=rss'
$ws.Range("I267").Value2 = '''$20'
$ws.Range("J267").Value2 = 'June 15, 2024 23:06 UTC'
$ws.Range("K267").Value2 = 'Video Editing'
$ws.Range("L267").Value2 = 'Video Editing'
$ws.Range("M267").Value2 = 'Singapore'

# --- Row 268 (new) ---
$ws.Range("A268").Value2 = '''51'
$ws.Range("B268").Value2 = 'YouTube Channel Manager - Upwork'
$ws.Range("C268").Value2 = 'https://www.upwork.com/jobs/YouTube-Channel-Manager_%7E01af099ff194680ce9?source=rss'
$v3 = @'
We are seeking a talented YouTube Channel Manager to oversee the growth and success of our channel. The ideal candidate will have a strong understanding of YouTube's algorithm and best practices for increasing engagement and subscribers. As the Channel Manager, you will be responsible for developing and executing a content strategy, optimizing video titles and tags, and analyzing performance metrics to make data-driven decisions. Additionally, you will collaborate with our creative team to ensure high-quality content production. 
  Skills needed:
  - Proficient in YouTube analytics and SEO
  - Excellent communication and organizational skills
  - Strong knowledge of social media marketing
  - Ability to analyze data and make strategic recommendations
  - Familiarity with video editing software
  This is a medium-sized project with a duration of 1 to 3 months. We are looking for an intermediate level expert who has prior experience managing successful YouTube channels.
Posted On
: June 15, 2024 23:05 UTC
Category
: Social Media Marketing
Skills
:YouTube Marketing,     YouTube Development,     YouTube,     Social Media Marketing,     Social Media Management    
Skills
:        YouTube Marketing,                     YouTube Development,                     YouTube,                     Social Media Marketing,                     Social Media Management            
Country
: United States
click to apply

'@
$ws.Range("D268").Value2 = $v3
$v4 = @'
We are seeking a talented YouTube Channel Manager to oversee the growth and success of our channel. The ideal candidate will have a strong understanding of YouTube&#039;s algorithm and best practices for increasing engagement and subscribers. As the Channel Manager, you will be responsible for developing and executing a content strategy, optimizing video titles and tags, and analyzing performance metrics to make data-driven decisions. Additionally, you will collaborate with our creative team to ensure high-quality content production. <br /><br />
&nbsp;&nbsp;Skills needed:<br />
&nbsp;&nbsp;- Proficient in YouTube analytics and SEO<br />
&nbsp;&nbsp;- Excellent communication and organizational skills<br />
&nbsp;&nbsp;- Strong knowledge of social media marketing<br />
&nbsp;&nbsp;- Ability to analyze data and make strategic recommendations<br />
&nbsp;&nbsp;- Familiarity with video editing software<br /><br />
&nbsp;&nbsp;This is a medium-sized project with a duration of 1 to 3 months. We are looking for an intermediate level expert who has prior experience managing successful YouTube channels.<br /><br /><br /><b>Posted On</b>: June 15, 2024 23:05 UTC<br /><b>Category</b>: Social Media Marketing<br /><b>Skills</b>:YouTube Marketing,     YouTube Development,     YouTube,     Social Media Marketing,     Social Media Management    
<br /><b>Skills</b>:        YouTube Marketing,                     YouTube Development,                     YouTube,                     Social Media Marketing,                     Social Media Management            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/YouTube-Channel-Manager_%7E01af099ff194680ce9?source=rss">click to apply</a>

'@
$ws.Range("E268").Value2 = $v4
$ws.Range("F268").Value2 = 'Sat, 15 Jun 2024 23:05:46 +0000'
$ws.Range("G268").Value2 = 'https://www.upwork.com/jobs/YouTube-Channel-Manager_%7E01af099ff194680ce9?source=rss'
$ws.Range("J268").Value2 = 'June 15, 2024 23:05 UTC'
$ws.Range("K268").Value2 = 'Social Media Marketing'
$ws.Range("L268").Value2 = 'YouTube Marketing,     YouTube Development,     YouTube,     Social Media Marketing,     Social Media Management'
$ws.Range("M268").Value2 = 'United States'

# --- Row 269 (new) ---
$ws.Range("A269").Value2 = '''51'
$ws.Range("B269").Value2 = 'Experienced YouTube Editor Needed - Upwork'
$ws.Range("C269").Value2 = 'https://www.upwork.com/jobs/Experienced-YouTube-Editor-Needed_%7E0115a94b4df68c84d4?source=rss'
$v5 = @'
We are looking for an experienced YouTube editor to join our team and help create engaging and professional video content for our channel. The ideal candidate should have a strong understanding of YouTube best practices and be able to edit videos that align with our brand and target audience. The main responsibilities will include editing raw footage, adding music and sound effects, creating eye-catching thumbnails, and implementing SEO strategies to optimize video visibility. The successful candidate should be proficient in video editing software such as Adobe Premiere Pro or Final Cut Pro. 
  Skills required:
  - Proficiency in video editing software (Adobe Premiere Pro, Final Cut Pro)
  - Strong understanding of YouTube best practices
  - Ability to create engaging and professional video content
  - Knowledge of SEO strategies for video optimization
Budget
: $80
Posted On
: June 15, 2024 22:55 UTC
Category
: Video Editing
Skills
:Video Editing,     Adobe Premiere Pro,     Audio Editing,     Video Post-Editing,     Adobe Photoshop    
Skills
:        Video Editing,                     Adobe Premiere Pro,                     Audio Editing,                     Video Post-Editing,                     Adobe Photoshop            
Country
: United States
click to apply

'@
$ws.Range("D269").Value2 = $v5
$v6 = @'
We are looking for an experienced YouTube editor to join our team and help create engaging and professional video content for our channel. The ideal candidate should have a strong understanding of YouTube best practices and be able to edit videos that align with our brand and target audience. The main responsibilities will include editing raw footage, adding music and sound effects, creating eye-catching thumbnails, and implementing SEO strategies to optimize video visibility. The successful candidate should be proficient in video editing software such as Adobe Premiere Pro or Final Cut Pro. <br />
&nbsp;&nbsp;Skills required:<br />
&nbsp;&nbsp;- Proficiency in video editing software (Adobe Premiere Pro, Final Cut Pro)<br />
&nbsp;&nbsp;- Strong understanding of YouTube best practices<br />
&nbsp;&nbsp;- Ability to create engaging and professional video content<br />
&nbsp;&nbsp;- Knowledge of SEO strategies for video optimization<br /><br /><b>Budget</b>: $80
<br /><b>Posted On</b>: June 15, 2024 22:55 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Video Editing,     Adobe Premiere Pro,     Audio Editing,     Video Post-Editing,     Adobe Photoshop    
<br /><b>Skills</b>:        Video Editing,                     Adobe Premiere Pro,                     Audio Editing,                     Video Post-Editing,                     Adobe Photoshop            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/Experienced-YouTube-Editor-Needed_%7E0115a94b4df68c84d4?source=rss">click to apply</a>

'@
$ws.Range("E269").Value2 = $v6
$ws.Range("F269").Value2 = 'Sat, 15 Jun 2024 22:55:32 +0000'
$ws.Range("G269").Value2 = 'https://www.upwork.com/jobs/Experienced-YouTube-Editor-Needed_%7E0115a94b4df68c84d4?source=rss'
$ws.Range("I269").Value2 = '''$80'
$ws.Range("J269").Value2 = 'June 15, 2024 22:55 UTC'
$ws.Range("K269").Value2 = 'Video Editing'
$ws.Range("L269").Value2 = 'Video Editing,     Adobe Premiere Pro,     Audio Editing,     Video Post-Editing,     Adobe Photoshop'
$ws.Range("M269").Value2 = 'United States'

# --- Row 270 (new) ---
$ws.Range("A270").Value2 = '''51'
$ws.Range("B270").Value2 = 'Quick Job -  Create A Zoom Frame(Border) For Me With My Logo - Upwork'
$ws.Range("C270").Value2 = 'https://www.upwork.com/jobs/Quick-Job-Create-Zoom-Frame-Border-For-With-Logo_%7E019b70b7d9c3ff466c?source=rss'
$v7 = @'
I need this really quickly
I'm looking for someone who can quickly create a Zoom frame(Border) for me similar to the one that I am attaching here:
I want it to have:
#1 My logo (On the left side)
#2- Text &quot;Get Your Ticket Today!&quot;(On the right side and above my logo)
#3- My website address: www.KCFLive.com (On the right side and large)
#4- An arrow pointing to the website
#5- I want you to use the follow colors (Black and gold, the same gold on my logo)
#6- I want it to be thin, so about the same size as the example that I sent you
Budget
: $15
Posted On
: June 15, 2024 22:45 UTC
Category
: Graphic Design
Skills
:Graphic Design,     Logo Design,     Adobe Photoshop,     Adobe Illustrator,     Illustration,     Adobe After Effects,     Web Design,     Video Editing,     Zoom Video Conferencing    
Skills
:        Graphic Design,                     Logo Design,                     Adobe Photoshop,                     Adobe Illustrator,                     Illustration,                     Adobe After Effects,                     Web Design,                     Video Editing,                     Zoom Video Conferencing            
Country
: United States
click to apply

'@
$ws.Range("D270").Value2 = $v7
$v8 = @'
I need this really quickly<br /><br />
I&#039;m looking for someone who can quickly create a Zoom frame(Border) for me similar to the one that I am attaching here:<br /><br />
I want it to have:<br /><br />
#1 My logo (On the left side)<br />
#2- Text &amp;quot;Get Your Ticket Today!&amp;quot;(On the right side and above my logo)<br />
#3- My website address: www.KCFLive.com (On the right side and large)<br />
#4- An arrow pointing to the website<br />
#5- I want you to use the follow colors (Black and gold, the same gold on my logo)<br />
#6- I want it to be thin, so about the same size as the example that I sent you<br /><br /><br /><b>Budget</b>: $15
<br /><b>Posted On</b>: June 15, 2024 22:45 UTC<br /><b>Category</b>: Graphic Design<br /><b>Skills</b>:Graphic Design,     Logo Design,     Adobe Photoshop,     Adobe Illustrator,     Illustration,     Adobe After Effects,     Web Design,     Video Editing,     Zoom Video Conferencing    
<br /><b>Skills</b>:        Graphic Design,                     Logo Design,                     Adobe Photoshop,                     Adobe Illustrator,                     Illustration,                     Adobe After Effects,                     Web Design,                     Video Editing,                     Zoom Video Conferencing            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/Quick-Job-Create-Zoom-Frame-Border-For-With-Logo_%7E019b70b7d9c3ff466c?source=rss">click to apply</a>

'@
$ws.Range("E270").Value2 = $v8
$ws.Range("F270").Value2 = 'Sat, 15 Jun 2024 22:45:33 +0000'
$ws.Range("G270").Value2 = 'https://www.upwork.com/jobs/Quick-Job-Create-Zoom-Frame-Border-For-With-Logo_%7E019b70b7d9c3ff466c?source=rss'
$ws.Range("I270").Value2 = '''$15'
$ws.Range("J270").Value2 = 'June 15, 2024 22:45 UTC'
$ws.Range("K270").Value2 = 'Graphic Design'
$ws.Range("L270").Value2 = 'Graphic Design,     Logo Design,     Adobe Photoshop,     Adobe Illustrator,     Illustration,     Adobe After Effects,     Web Design,     Video Editing,     Zoom Video Conferencing'
$ws.Range("M270").Value2 = 'United States'

# --- Row 271 (new) ---
$ws.Range("A271").Value2 = '''51'
$ws.Range("B271").Value2 = 'I need a white square overlay on a 3 minute video and sound removed and replaced with music. - Upwork'
$ws.Range("C271").Value2 = 'https://www.upwork.com/jobs/need-white-square-overlay-minute-video-and-sound-removed-and-replaced-with-music_%7E01ec9810f5a110b541?source=rss'
$v9 = @'
I need a white square overlay in the top right had corner to cover up the presenter.
After the guy stops talking in the beginning, remove the sound and add back good sound.
Budget
: $5
Posted On
: June 15, 2024 22:40 UTC
Category
: Video Editing
Skills
:Video Editing    
Skills
:        Video Editing            
Country
: United States
click to apply

'@
$ws.Range("D271").Value2 = $v9
$v10 = @'
I need a white square overlay in the top right had corner to cover up the presenter.<br /><br />
After the guy stops talking in the beginning, remove the sound and add back good sound.<br /><br /><b>Budget</b>: $5
<br /><b>Posted On</b>: June 15, 2024 22:40 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Video Editing    
<br /><b>Skills</b>:        Video Editing            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/need-white-square-overlay-minute-video-and-sound-removed-and-replaced-with-music_%7E01ec9810f5a110b541?source=rss">click to apply</a>

'@
$ws.Range("E271").Value2 = $v10
$ws.Range("F271").Value2 = 'Sat, 15 Jun 2024 22:40:51 +0000'
$ws.Range("G271").Value2 = 'https://www.upwork.com/jobs/need-white-square-overlay-minute-video-and-sound-removed-and-replaced-with-music_%7E01ec9810f5a110b541?source=rss'
$ws.Range("I271").Value2 = '''$5'
$ws.Range("J271").Value2 = 'June 15, 2024 22:40 UTC'
$ws.Range("K271").Value2 = 'Video Editing'
$ws.Range("L271").Value2 = 'Video Editing'
$ws.Range("M271").Value2 = 'United States'
